# "Merge Unmerge Insert Col & Row" -- final state:
#   B1 = long "Hello..." string, merged B1:E1
#   B2 = "World!",              merged B2:C2
#   B3 = blank (kept as an explicit empty cell)
#   B4 = "Foo"
#   B5 = "Bar"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Overwrite the numeric values with text.
$ws.Range("B1").Value = "Hello, this is a very very long string."
$ws.Range("B2").Value = "World!"

# Row 3 becomes blank, but keep the cell itself alive (re-assert the sheet's
# own default vertical alignment so no new style is introduced) instead of
# letting it disappear completely.
$ws.Range("B3").ClearContents()
$ws.Range("B3").VerticalAlignment = -4107

# Write the new values for rows 5 and 4 -- Bar first so it lands earlier in
# the shared-string table than Foo.
$ws.Range("B5").Value = "Bar"
$ws.Range("B4").Value = "Foo"

# Merge the header-ish rows.
$ws.Range("B1:E1").Merge()
$ws.Range("B2:C2").Merge()
